$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $lastRow = $ws.UsedRange.Rows.Count

    $ws.Cells.Item(2, 1).Value = 17
    $ws.Cells.Item(3, 1).Value = 18
    if ($lastRow -ge 4) {
        $ws.Cells.Item(4, 1).Value = 19
    }
}
